$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold numeric-looking text values (e.g. "303.01",
# "2.57%") that must be preserved as literal text, matching the original workbook's
# inline-string encoding. Force those cells to Text format first so Excel does not
# silently coerce them into numbers/percentages.
$textCells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9",
    "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "D15", "E15", "D16", "E16",
    "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23",
    "E23", "D24", "E24", "D25", "E25", "E26", "E27", "D39", "E39", "D40", "E40", "D41", "E41",
    "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "E47", "D48", "E48",
    "D49", "E49", "D50", "E50", "D51", "E51"
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated cell values (prices, % volumes, coin names and links) so the sheet
# matches the refreshed crypto symbol list.
$ws.Range("D2").Value = "303.01"
$ws.Range("E2").Value = "2.57%"
$ws.Range("D3").Value = "35.06"
$ws.Range("E3").Value = "12.65%"
$ws.Range("D4").Value = "5.151"
$ws.Range("E4").Value = "4.44%"
$ws.Range("E5").Value = "4.34%"
$ws.Range("D6").Value = "2.362"
$ws.Range("E6").Value = "5.80%"
$ws.Range("D7").Value = "8.047"
$ws.Range("E7").Value = "3.60%"
$ws.Range("D8").Value = "3.944"
$ws.Range("E8").Value = "5.33%"
$ws.Range("D9").Value = "0.9290"
$ws.Range("E9").Value = "1.71%"
$ws.Range("D10").Value = "0.1003"
$ws.Range("E10").Value = "10.85%"
$ws.Range("D11").Value = "0.1798"
$ws.Range("E11").Value = "4.87%"
$ws.Range("D12").Value = "0.08597"
$ws.Range("E12").Value = "3.79%"
$ws.Range("D13").Value = "0.03316"
$ws.Range("E13").Value = "5.68%"
$ws.Range("D14").Value = "0.09924"
$ws.Range("D15").Value = "0.001498"
$ws.Range("E15").Value = "-1.18%"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "0.04563"
$ws.Range("E16").Value = "0.26%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "0.005802"
$ws.Range("E17").Value = "1.02%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "3.462"
$ws.Range("E18").Value = "-1.07%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "2.166"
$ws.Range("E19").Value = "4.39%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "0.3359"
$ws.Range("E20").Value = "0.95%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "0.1333"
$ws.Range("E21").Value = "2.77%"
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").Value = "4.286"
$ws.Range("E22").Value = "7.65%"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "0.2302"
$ws.Range("E23").Value = "9.60%"
$ws.Range("D24").Value = "0.001213"
$ws.Range("E24").Value = "0.06%"
$ws.Range("D25").Value = "0.004370"
$ws.Range("E25").Value = "-5.37%"
$ws.Range("E26").Value = "-0.01%"
$ws.Range("E27").Value = "0.04%"
$ws.Range("D39").Value = "0.01790"
$ws.Range("E39").Value = "10.95%"
$ws.Range("D40").Value = "0.04795"
$ws.Range("E40").Value = "6.74%"
$ws.Range("D41").Value = "0.007794"
$ws.Range("E41").Value = "7.69%"
$ws.Range("D42").Value = "0.1414"
$ws.Range("E42").Value = "6.31%"
$ws.Range("D43").Value = "0.007333"
$ws.Range("E43").Value = "-18.52%"
$ws.Range("D44").Value = "0.002139"
$ws.Range("E44").Value = "8.61%"
$ws.Range("D45").Value = "0.009454"
$ws.Range("E45").Value = "3.78%"
$ws.Range("D46").Value = "0.00006109"
$ws.Range("E46").Value = "0.31%"
$ws.Range("E47").Value = "-0.01%"
$ws.Range("D48").Value = "3.060"
$ws.Range("E48").Value = "38.08%"
$ws.Range("D49").Value = "0.001999"
$ws.Range("E49").Value = "-0.03%"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").Value = "-0.01%"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").Value = "-0.01%"
